$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2024-03-02 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-03-03 Sunday", 2) | Out-Null

# Update the multiplication problems in the table.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "562×2="
$t.Cell(1, 2).Range.Text = "829×2="
$t.Cell(1, 3).Range.Text = "293×5="
$t.Cell(1, 4).Range.Text = "752×7="
$t.Cell(1, 5).Range.Text = "800×2="

$t.Cell(5, 1).Range.Text = "573×4="
$t.Cell(5, 2).Range.Text = "734×9="
$t.Cell(5, 3).Range.Text = "541×2="
$t.Cell(5, 4).Range.Text = "613×9="
$t.Cell(5, 5).Range.Text = "567×6="

$t.Cell(10, 1).Range.Text = "576×3="
$t.Cell(10, 2).Range.Text = "707×3="
$t.Cell(10, 3).Range.Text = "906×7="
$t.Cell(10, 4).Range.Text = "716×3="
$t.Cell(10, 5).Range.Text = "949×8="

$t.Cell(15, 1).Range.Text = "410×7="
$t.Cell(15, 2).Range.Text = "820×4="
$t.Cell(15, 3).Range.Text = "740×4="
$t.Cell(15, 4).Range.Text = "418×5="
$t.Cell(15, 5).Range.Text = "631×9="

$t.Cell(20, 1).Range.Text = "936×7="
$t.Cell(20, 2).Range.Text = "401×7="
$t.Cell(20, 3).Range.Text = "891×5="
$t.Cell(20, 4).Range.Text = "871×2="
$t.Cell(20, 5).Range.Text = "119×6="
